# Add a new "2021" column (R) of data to the SDG 1.3.1 indicator sheet,
# mirroring the existing year columns (D..Q = 2007..2020).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

function Copy-FormatTo {
    param(
        [string]$SourceAddress,
        [string]$DestAddress
    )
    $ws.Range($SourceAddress).Copy() | Out-Null
    $ws.Range($DestAddress).PasteSpecial($xlPasteFormats) | Out-Null
}

# Row 2 - bottom border filler cell under the year header, no value, same
# style as the one to its left (Q2).
Copy-FormatTo "Q2" "R2"

# Row 3 - year header "2021", same style as the other year headers.
Copy-FormatTo "D3" "R3"
$ws.Range("R3").Value = 2021

# Row 4 - headline indicator row; existing cells use style s="27" (General
# number format). The new column keeps that look but with a "0.0" number
# format applied on top of it.
Copy-FormatTo "Q4" "R4"
$ws.Range("R4").Value = 18
$ws.Range("R4").NumberFormat = "0.0"

# Rows 5-12 - data rows sharing the "0.0" style used by columns D..O
# (s="25") rather than Q's slightly different xf (s="28").
Copy-FormatTo "D5" "R5"
$ws.Range("R5").Value = 1.7480265877296817

Copy-FormatTo "D6" "R6"
$ws.Range("R6").Value = 4.1112601249414027

Copy-FormatTo "D7" "R7"
$ws.Range("R7").Value = 1.5225742120245318

Copy-FormatTo "D8" "R8"
$ws.Range("R8").Value = 1.2326518235454269

Copy-FormatTo "D9" "R9"
$ws.Range("R9").Value = 4.0865392096984241

Copy-FormatTo "D10" "R10"
$ws.Range("R10").Value = 0.84876624403485645

Copy-FormatTo "D11" "R11"
$ws.Range("R11").Value = 2.1456657699653627

Copy-FormatTo "D12" "R12"
$ws.Range("R12").Value = 1.8214779402142154

# Row 13 - bottom (bordered) data row, style s="26".
Copy-FormatTo "D13" "R13"
$ws.Range("R13").Value = 0.51989507542472779

# Leave the same selection the authored workbook ends up with.
$ws.Range("R24:R25").Select() | Out-Null
